# ADD results from server
# Update fix-cost result values (row 2) on each year sheet with the latest
# server-computed figures.

$wb = $excel.ActiveWorkbook

$updates = @{
    "2025" = @{
        "B2" = 1037.265132737054
        "E2" = 28926.05393052954
        "G2" = 8095.925712661834
        "I2" = 16171.06685703679
        "L2" = 48492.22142001599
        "M2" = 10595.37713982
        "N2" = 7071.74531360843
        "O2" = 6993.890772562212
    }
    "2030" = @{
        "A2" = 0
        "B2" = 4157.588990853394
        "E2" = 45991.90904307188
        "G2" = 8095.925712661834
        "I2" = 37079.12819938764
        "L2" = 54844.03303316472
        "M2" = 17449.04999683176
        "N2" = 9024.733389685653
        "O2" = 9724.258249348202
    }
    "2035" = @{
        "A2" = 2754.31755456332
        "B2" = 6368.910634126893
        "E2" = 57457.45307013817
        "G2" = 8095.925712661834
        "I2" = 52465.73681402855
        "L2" = 54844.03303316472
        "M2" = 21912.87293902603
        "N2" = 13034.3101291405
        "O2" = 12860.17168993684
    }
    "2040" = @{
        "A2" = 2754.31755456332
        "B2" = 6368.910634126893
        "E2" = 57457.45307013817
        "G2" = 8095.925712661834
        "I2" = 52465.73681402855
        "L2" = 54844.03303316472
        "M2" = 21912.87293902603
        "N2" = 13149.40378733627
        "O2" = 12860.17168993684
    }
    "2045" = @{
        "A2" = 5713.151062849596
        "B2" = 6368.910634126893
        "E2" = 57457.45307013817
        "G2" = 8095.925712661834
        "I2" = 52465.73681402855
        "L2" = 54844.03303316472
        "M2" = 21912.87293902603
        "N2" = 13598.62114148921
        "O2" = 14937.1305943757
    }
    "2050" = @{
        "A2" = 5713.151062849596
        "B2" = 6368.910634126893
        "E2" = 57457.45307013817
        "G2" = 8095.925712661834
        "I2" = 52465.73681402855
        "L2" = 54844.03303316472
        "M2" = 21912.87293902603
        "N2" = 13598.62114148921
        "O2" = 14937.1305943757
    }
}

foreach ($sheetName in $updates.Keys) {
    $sheetNameStr = [string]$sheetName
    $ws = $wb.Worksheets.Item($sheetNameStr)
    $cellValues = $updates[$sheetName]
    foreach ($cellRef in $cellValues.Keys) {
        $ws.Range([string]$cellRef).Value = $cellValues[$cellRef]
    }
}
